# Updates the LinkedIn carousel draft: swaps the Elcogen/Saur Energy
# article content for the Insolation Energy / Bondada Engineering MOU story
# across all 6 slides (title + two body bullet points per slide).

$p = $ppt.ActivePresentation

$newTitle = "Insolation Energy Subsidiary Signs MOU With Bondada Engineering For Solar Modules - scanx.trade"

# Per-slide replacement text for the 2nd and 3rd paragraphs of the body
# placeholder (the 1st paragraph, "Image missing", is left untouched).
$bodyReplacements = @(
    @("Insolation Energy has signed a Memorandum of Understanding (MOU) with Bondada Engineering.", "The MOU focuses on the development of solar modules."),
    @("The partnership aims to enhance solar energy solutions.", "Bondada Engineering is recognized for its expertise in engineering services."),
    @("The collaboration is expected to contribute to renewable energy initiatives.", "Details regarding the specific solar module technology were not disclosed."),
    @("The MOU signifies a strategic partnership in the solar energy sector.", "Insolation Energy's subsidiary is focused on expanding its market presence."),
    @("The agreement is part of a broader trend towards sustainable energy solutions.", "Both companies aim to leverage their strengths in the renewable energy market."),
    @("The MOU could lead to future projects and collaborations in solar technology.", "This partnership aligns with global efforts to increase solar energy adoption.")
)

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)

    # Shape 1: Title. Clear first, then set (via the single-paragraph
    # sub-range, same trick as the body below) so the runtime doesn't try to
    # diff-preserve a run split against the old text (e.g. shared trailing
    # punctuation) and instead emits a single clean <a:r> with no extra
    # <a:rPr>.
    $titleShape = $slide.Shapes.Item(1)
    $titleTextRange = $titleShape.TextFrame.TextRange
    $titleTextRange.Paragraphs(1, 1).Text = ""
    $titleTextRange.Paragraphs(1, 1).Text = $newTitle

    # Shape 2: Content placeholder - update paragraphs 2 and 3 in place,
    # leaving paragraph 1 ("Image missing") untouched.
    $bodyShape = $slide.Shapes.Item(2)
    $bodyTextRange = $bodyShape.TextFrame.TextRange

    $pair = $bodyReplacements[$i - 1]

    $bodyTextRange.Paragraphs(2, 1).Text = ""
    $bodyTextRange.Paragraphs(2, 1).Text = $pair[0]

    $bodyTextRange.Paragraphs(3, 1).Text = ""
    $bodyTextRange.Paragraphs(3, 1).Text = $pair[1]
}
